$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 54
$ws1.Range("F3").Value = 3175
$ws1.Range("F5").Value = 2209
$ws1.Range("F6").Value = 330
$ws1.Range("F8").Value = 1059
$ws1.Range("F9").Value = 1021
$ws1.Range("F10").Value = 245
$ws1.Range("F11").Value = 465
$ws1.Range("F12").Value = 1160
$ws1.Range("F14").Value = 77
$ws1.Range("F16").Value = 7856
$ws1.Range("F17").Value = 344
$ws1.Range("F18").Value = 2467
$ws1.Range("F19").Value = 218
$ws1.Range("F20").Value = 232
$ws1.Range("F23").Value = 539
$ws1.Range("F27").Value = 1520
$ws1.Range("F28").Value = 6
$ws1.Range("F30").Value = 1662
$ws1.Range("F32").Value = 1909
$ws1.Range("F34").Value = 45
$ws1.Range("F35").Value = 168
$ws1.Range("F36").Value = 276
$ws1.Range("F37").Value = 43
$ws1.Range("F38").Value = 177
$ws1.Range("F39").Value = 347
$ws1.Range("F41").Value = 219

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 54
$ws4.Range("F5").Value = 3175
$ws4.Range("F7").Value = 2209
$ws4.Range("F8").Value = 330
$ws4.Range("F10").Value = 1059
$ws4.Range("F12").Value = 1021
$ws4.Range("F13").Value = 245
$ws4.Range("F14").Value = 465
$ws4.Range("F15").Value = 1160
$ws4.Range("F17").Value = 77
$ws4.Range("F19").Value = 7856
$ws4.Range("F20").Value = 344
$ws4.Range("F21").Value = 2467
$ws4.Range("F23").Value = 218
$ws4.Range("F24").Value = 232
$ws4.Range("F27").Value = 539
$ws4.Range("F31").Value = 1520
$ws4.Range("F32").Value = 6
$ws4.Range("F34").Value = 1662
$ws4.Range("F36").Value = 1909
$ws4.Range("F38").Value = 45
$ws4.Range("F39").Value = 168
$ws4.Range("F40").Value = 276
$ws4.Range("F41").Value = 43
$ws4.Range("F42").Value = 177
$ws4.Range("F43").Value = 347
$ws4.Range("F48").Value = 219

$wb.Save()
